$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42 (shifts existing rows 42:91 down to 43:92)
$ws.Rows.Item(42).Insert()

# New row 42 mirrors the old row 42 data, but with an updated date (D) and volume (J)
$ws.Cells.Item(42, 1).Value = 9
$ws.Cells.Item(42, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(42, 3).Value = "Metropolitana"
$ws.Cells.Item(42, 4).Value = 44671
$ws.Cells.Item(42, 5).Value = 13
$ws.Cells.Item(42, 6).Value = 100112005
$ws.Cells.Item(42, 7).Value = "Puerro"
$ws.Cells.Item(42, 8).Value = "Sin especificar"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 160
$ws.Cells.Item(42, 11).Value = 8000
$ws.Cells.Item(42, 12).Value = 8000
$ws.Cells.Item(42, 13).Value = 8000
$ws.Cells.Item(42, 14).Value = "$/paquete 20 unidades"
$ws.Cells.Item(42, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(42, 16).Value = 400
$ws.Cells.Item(42, 17).Value = 20
$ws.Cells.Item(42, 18).Value = "Hortaliza"
